$d = $word.ActiveDocument

# "Christopher Thacker" -> "Chris T."
$d.Content.Find.Execute("Christopher Thacker", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Chris T.", 2)

# "Ioannis Batsois" (was split across three runs with a spell-check proof
# mark around each name) -> single run "Ioannis Batsois"
$d.Content.Find.Execute("Ioannis Batsois", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ioannis Batsois", 2)
